# implemented CR, CP, F, AR metrics and generated consolidated report
#
# Summary sheet: add a 10th run (run 10) as a new data row, add a new
# "execution_time_seconds" column (J) with the new run's execution time,
# and recompute the TOTALS row. Question Failure Rates sheet: each
# question now reflects 10 total runs instead of 9 (updating passed /
# failed counts and the failure_rate text accordingly).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Summary")
$ws2 = $wb.Worksheets.Item("Question Failure Rates")

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------

# Insert a fresh row 11 (pushes the old TOTALS row from 11 down to 12)
# so the new run's data lands right after run 9 and TOTALS stays last.
$ws1.Rows.Item(11).Insert()

# New column J: execution_time_seconds (copy header cell formatting from I1)
$ws1.Range("I1").Copy($ws1.Range("J1"))
$ws1.Range("J1").Value = "execution_time_seconds"

# New row 11: results for run 10
$ws1.Range("A11").Value = 10
$ws1.Range("B11").Value = 35
$ws1.Range("C11").Value = 30
$ws1.Range("D11").Value = 5
$ws1.Range("E11").Value = 0
$ws1.Range("F11").Value = "'85.7%"
$ws1.Range("G11").Value = 0.6438
$ws1.Range("H11").Value = "Q, R, S, Q28, Q29"
$ws1.Range("I11").Value = "-"
$ws1.Range("J11").Value = 141.94

# TOTALS row (now row 12): 10 runs, recomputed pass rate / avg score,
# plus the total execution time in the new column.
$ws1.Range("B12").Value = "10 runs"
$ws1.Range("F12").Value = "'76.6%"
$ws1.Range("G12").Value = 0.5856
$ws1.Range("J12").Value = 141.94

# ---------------------------------------------------------------------
# Question Failure Rates sheet: every question now has 10 total_runs
# ---------------------------------------------------------------------

$ws2.Range("C2").Value = 10
$ws2.Range("D2").Value = 10

$ws2.Range("C3").Value = 10
$ws2.Range("D3").Value = 8
$ws2.Range("G3").Value = "'20.0%"

$ws2.Range("C4").Value = 10
$ws2.Range("D4").Value = 10

$ws2.Range("C5").Value = 10
$ws2.Range("D5").Value = 3
$ws2.Range("G5").Value = "'60.0%"

$ws2.Range("C6").Value = 10
$ws2.Range("D6").Value = 8
$ws2.Range("G6").Value = "'20.0%"

$ws2.Range("C7").Value = 10
$ws2.Range("D7").Value = 8
$ws2.Range("G7").Value = "'20.0%"

$ws2.Range("C8").Value = 10
$ws2.Range("D8").Value = 10

$ws2.Range("C9").Value = 10
$ws2.Range("D9").Value = 10

$ws2.Range("C10").Value = 10
$ws2.Range("D10").Value = 10

$ws2.Range("C11").Value = 10
$ws2.Range("D11").Value = 5
$ws2.Range("G11").Value = "'50.0%"

$ws2.Range("C12").Value = 10
$ws2.Range("D12").Value = 10

$ws2.Range("C13").Value = 10
$ws2.Range("D13").Value = 8
$ws2.Range("G13").Value = "'20.0%"

$ws2.Range("C14").Value = 10
$ws2.Range("D14").Value = 9
$ws2.Range("G14").Value = "'10.0%"

$ws2.Range("C15").Value = 10
$ws2.Range("D15").Value = 9
$ws2.Range("G15").Value = "'10.0%"

$ws2.Range("C16").Value = 10
$ws2.Range("D16").Value = 10

$ws2.Range("C17").Value = 10
$ws2.Range("D17").Value = 10

$ws2.Range("C18").Value = 10
$ws2.Range("E18").Value = 5
$ws2.Range("G18").Value = "'50.0%"

$ws2.Range("C19").Value = 10
$ws2.Range("E19").Value = 10

$ws2.Range("C20").Value = 10
$ws2.Range("E20").Value = 10

$ws2.Range("C21").Value = 10
$ws2.Range("D21").Value = 9
$ws2.Range("G21").Value = "'10.0%"

$ws2.Range("C22").Value = 10
$ws2.Range("D22").Value = 9
$ws2.Range("G22").Value = "'10.0%"

$ws2.Range("C23").Value = 10
$ws2.Range("D23").Value = 8
$ws2.Range("G23").Value = "'20.0%"

$ws2.Range("C24").Value = 10
$ws2.Range("D24").Value = 8
$ws2.Range("G24").Value = "'20.0%"

$ws2.Range("C25").Value = 10
$ws2.Range("D25").Value = 8
$ws2.Range("G25").Value = "'20.0%"

$ws2.Range("C26").Value = 10
$ws2.Range("D26").Value = 8
$ws2.Range("G26").Value = "'20.0%"

$ws2.Range("C27").Value = 10
$ws2.Range("D27").Value = 8
$ws2.Range("G27").Value = "'20.0%"

$ws2.Range("C28").Value = 10
$ws2.Range("D28").Value = 9
$ws2.Range("G28").Value = "'10.0%"

$ws2.Range("C29").Value = 10
$ws2.Range("E29").Value = 2
$ws2.Range("G29").Value = "'20.0%"

$ws2.Range("C30").Value = 10
$ws2.Range("E30").Value = 10

$ws2.Range("C31").Value = 10
$ws2.Range("D31").Value = 10

$ws2.Range("C32").Value = 10
$ws2.Range("D32").Value = 10

$ws2.Range("C33").Value = 10
$ws2.Range("D33").Value = 10

$ws2.Range("C34").Value = 10
$ws2.Range("D34").Value = 6
$ws2.Range("G34").Value = "'40.0%"

$ws2.Range("C35").Value = 10
$ws2.Range("D35").Value = 8
$ws2.Range("G35").Value = "'20.0%"

$ws2.Range("C36").Value = 10
$ws2.Range("D36").Value = 9
$ws2.Range("G36").Value = "'10.0%"
